$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark additional wishlist items as "Reserved / bought" by writing "Y" in column E
$ws.Range("E2").Value = "Y"
$ws.Range("E4").Value = "Y"
$ws.Range("E6").Value = "Y"
$ws.Range("E11").Value = "Y"
$ws.Range("E12").Value = "Y"
$ws.Range("E13").Value = "Y"

# Match the final cursor/selection position recorded in the workbook
$ws.Range("E13").Select()
